# Update "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column after "Week" (new column B)
#  - shorten the Week labels (W01 -> W1, etc.)
#  - shift/update the forecast numbers to their corrected values
#  - make the is_holiday_week column a boolean column

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1. Insert a new column B ("Week_Start_Date"), shifting ASIN..is_holiday_week right by one.
$ws.Columns.Item(2).Insert()

# Make sure the new date column stores plain text (not auto-converted Excel dates).
$ws.Columns.Item(2).NumberFormat = "@"

# 2. Header row
$ws.Range("B1").Value = "Week_Start_Date"

# 3. Row data: Week, Week_Start_Date, MyForecast, AmazonMean, P70, P80, P90
$data = @(
    @("W1",  "2025-01-05", 2, 2, 2, 3, 6),
    @("W2",  "2025-01-12", 2, 2, 2, 3, 6),
    @("W3",  "2025-01-19", 2, 2, 2, 3, 6),
    @("W4",  "2025-01-26", 2, 2, 2, 3, 6),
    @("W5",  "2025-02-02", 2, 2, 2, 3, 6),
    @("W6",  "2025-02-09", 2, 2, 2, 3, 7),
    @("W7",  "2025-02-16", 2, 3, 2, 4, 9),
    @("W8",  "2025-02-23", 2, 2, 2, 3, 6),
    @("W9",  "2025-03-02", 2, 2, 2, 4, 7),
    @("W10", "2025-03-09", 2, 3, 2, 4, 9),
    @("W11", "2025-03-16", 2, 3, 2, 5, 10),
    @("W12", "2025-03-23", 3, 4, 3, 6, 13),
    @("W13", "2025-03-30", 3, 4, 3, 6, 13),
    @("W14", "2025-04-06", 3, 4, 3, 6, 13),
    @("W15", "2025-04-13", 3, 4, 3, 7, 14),
    @("W16", "2025-04-20", 3, 4, 3, 6, 13)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]        # A: Week
    $ws.Cells.Item($r, 2).Value = $row[1]        # B: Week_Start_Date
    $ws.Cells.Item($r, 4).Value = $row[2]        # D: MyForecast
    $ws.Cells.Item($r, 5).Value = $row[3]        # E: Amazon Mean Forecast
    $ws.Cells.Item($r, 6).Value = $row[4]        # F: Amazon P70 Forecast
    $ws.Cells.Item($r, 7).Value = $row[5]        # G: Amazon P80 Forecast
    $ws.Cells.Item($r, 8).Value = $row[6]        # H: Amazon P90 Forecast
    $ws.Cells.Item($r, 10).Value = $false        # J: is_holiday_week (boolean)
    $r = $r + 1
}
